# Apply the commit: "Added power sequencing for ADC. Split IMU into 2 boards."
#
# Summary of the change:
#   - The old 3.3V LM1117DT-1.8 linear regulator part (row 6) is replaced by a
#     TPS-1.8 switching regulator with new current/voltage numbers.
#   - Four new summary rows are added (21-24) that split the previous single
#     IMU board totals into an "ANALOG BOARD" and a "DIGITAL BOARD", each with
#     their own 3.3V and 5V current rollups, to reflect the 2-board split
#     mentioned in the commit message ("power sequencing for ADC").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary rows first (21-24), so new shared strings are appended in
#     the same order the original author typed them (board totals before the
#     regulator rename) ---
$ws.Range("A21").Value() = "ANALOG BOARD 3.3V Current"
$ws.Range("B21").Formula() = "=SUM(F2,F8,F6)"

$ws.Range("A22").Value() = "ANALOG BOARD 5V Current"
$ws.Range("B22").Formula() = "=SUM(F5,F4,F3,F7)"

$ws.Range("A23").Value() = "DIGITAL BOARD 3.3V Current"
$ws.Range("B23").Formula() = "=SUM(F11,F10,F9)"

$ws.Range("A24").Value() = "DIGITAL BOARD 5V Current"
$ws.Range("B24").Value() = 0

# --- Replace the old LM1117DT-1.8 part (row 6) with the new TPS-1.8 part ---
$ws.Range("A6").Value() = "TPS-1.8"
$ws.Range("B6").Value() = 20
$ws.Range("C6").Value() = 3.3

# --- Update the selection to match the final state of the workbook ---
$ws.Range("B25").Select()
